$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fill in Priority and Type for row 21 (TC for "Select a product and add it to the shopping cart.")
$ws.Range("D21").Value = "High"
$ws.Range("E21").Value = "Sanity, Regression"

# Update the selected cell to A22, matching the saved selection state
$ws.Range("A22").Select()
